$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 136, shifting the existing rows 136-168 down to 137-169.
$ws.Rows("136").Insert()

# Populate the newly inserted row 136 with its data.
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 44889
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112032
$ws.Range("G136").Value = "Zapallo italiano"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 220
$ws.Range("K136").Value = 5500
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 5727
$ws.Range("N136").Value = "$/caja 50 unidades"
$ws.Range("O136").Value = "Región de O'Higgins"
$ws.Range("P136").Value = 115
$ws.Range("Q136").Value = 50
$ws.Range("R136").Value = "Hortaliza"
